$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs target) - D2 unchanged text "ECs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.06815433333333333
$ws.Range("H2").Value2 = 0.204463
$ws.Range("M2").Value2 = 8.799974000000001
$ws.Range("N2").Value2 = 26.399922
$ws.Range("O2").Value2 = 0.1502220177021807
$ws.Range("P2").Value2 = 0.1502220177021807
$ws.Range("Q2").Value2 = 0.5997563613206667
$ws.Range("R2").Value2 = 5.397807251886
$ws.Range("S2").Value2 = 0.1502220177021807
$ws.Range("T2").Value2 = 0.1502220177021807

# Row 3 (FAPs target) - D3 unchanged text "FAPs"
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.06815433333333333
$ws.Range("H3").Value2 = 0.204463
$ws.Range("O3").Value2 = 0.4499951903206205
$ws.Range("P3").Value2 = 0.4499951903206205
$ws.Range("Q3").Value2 = 1.796590686816333
$ws.Range("R3").Value2 = 16.169316181347
$ws.Range("S3").Value2 = 0.4499951903206205
$ws.Range("T3").Value2 = 0.4499951903206205

# Row 4 (MuSCs target) - D4 text stays "MuSCs"
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 0.06815433333333333
$ws.Range("H4").Value2 = 0.204463
$ws.Range("M4").Value2 = 23.07309566666667
$ws.Range("N4").Value2 = 69.21928700000001
$ws.Range("O4").Value2 = 0.3938746848208995
$ws.Range("P4").Value2 = 0.3938746848208995
$ws.Range("Q4").Value2 = 1.572531453097889
$ws.Range("R4").Value2 = 14.152783077881
$ws.Range("S4").Value2 = 0.3938746848208995
$ws.Range("T4").Value2 = 0.3938746848208995

# Row 5 (Resolving-Mac target) - D5 text stays "Resolving-Mac"
$ws.Range("D5").Value2 = "Resolving-Mac"
$ws.Range("E5").Value2 = 2
$ws.Range("F5").Value2 = 0.6666666666666666
$ws.Range("G5").Value2 = 0.06815433333333333
$ws.Range("H5").Value2 = 0.204463
$ws.Range("M5").Value2 = 0.3460956666666666
$ws.Range("N5").Value2 = 1.038287
$ws.Range("O5").Value2 = 0.005908107156299329
$ws.Range("P5").Value2 = 0.00590810715629933
$ws.Range("Q5").Value2 = 0.02358791943122222
$ws.Range("R5").Value2 = 0.212291274881
$ws.Range("S5").Value2 = 0.005908107156299329
$ws.Range("T5").Value2 = 0.00590810715629933
